# Swap the roles of the two sheets:
#   - the sheet currently named "hotel_info" becomes "review_info" and is
#     emptied down to just the review_info header row (25 columns, A:Y)
#   - the sheet currently named "review_info" becomes "hotel_info" and is
#     populated with the hotel_info header row + one data row, with a new
#     "State" column inserted between Hotel_Name and City.

$wb = $excel.ActiveWorkbook

$hotelSheet  = $wb.Worksheets.Item("hotel_info")
$reviewSheet = $wb.Worksheets.Item("review_info")

# --- Rename sheets (keep tab order/position, just swap the names) ---
# use a temporary name to avoid a collision while swapping
$hotelSheet.Name = "tmp_swap_name"
$reviewSheet.Name = "hotel_info"
$hotelSheet.Name = "review_info"

# Now:
#   $hotelSheet  (was tab 1, rId1) -> Name = "review_info"
#   $reviewSheet (was tab 2, rId2) -> Name = "hotel_info"

# --- Rebuild the "review_info" sheet content (was "hotel_info") ---
$newReviewSheet = $wb.Worksheets.Item("review_info")
$newReviewSheet.Cells.Clear()

$reviewHeaders = @(
    "STR",
    "reviewer_ID",
    "reviewer_name",
    "Review_ID",
    "Date_of_scraping",
    "ReviewURL",
    "Tripadvisor_gcode",
    "Tripadvisor_dcode",
    "Tripadvisor_rcode",
    "review_date",
    "review_title",
    "review_content",
    "review_rating",
    "trip_month",
    "trip_purpose",
    "value",
    "rooms",
    "Location",
    "Cleanliness",
    "Sleep Quality",
    "Service",
    "Picture(yes=1)",
    "respondent",
    "response_date",
    "response_text"
)

for ($i = 0; $i -lt $reviewHeaders.Count; $i++) {
    $newReviewSheet.Cells.Item(1, $i + 1).Value = $reviewHeaders[$i]
}

# --- Rebuild the "hotel_info" sheet content (was "review_info") ---
$newHotelSheet = $wb.Worksheets.Item("hotel_info")
$newHotelSheet.Cells.Clear()

$hotelHeaders = @(
    "STR",
    "Hotel_Name",
    "State",
    "City",
    "Zip",
    "TA_ReviewURL",
    "Tripadvisor_Hotel_Name",
    "English_Reviews_num",
    "Local_Rank",
    "Total_Reviews_num"
)

for ($i = 0; $i -lt $hotelHeaders.Count; $i++) {
    $newHotelSheet.Cells.Item(1, $i + 1).Value = $hotelHeaders[$i]
}

# English_Reviews_num / Local_Rank / Total_Reviews_num are stored as TEXT
# (shared-string) values in the source data, not numbers - force the
# "Text" number format on those cells before assigning so the numeric-
# looking strings ("1027", "91", "1047") aren't auto-coerced to numbers.
$newHotelSheet.Range("H2:J2").NumberFormat = "@"

$newHotelSheet.Cells.Item(2, 1).Value = 42234
$newHotelSheet.Cells.Item(2, 2).Value = "Hilton New Orleans St Charles Avenue"
$newHotelSheet.Cells.Item(2, 3).Value = "Louisiana"
$newHotelSheet.Cells.Item(2, 4).Value = "New Orleans"
$newHotelSheet.Cells.Item(2, 5).Value = 70130
$newHotelSheet.Cells.Item(2, 6).Value = "https://www.tripadvisor.com/Hotel_Review-g60864-d638900-Reviews-Hilton_New_Orleans_St_Charles_Avenue-New_Orleans_Louisiana.html"
$newHotelSheet.Cells.Item(2, 7).Value = "Hilton New Orleans/St. Charles Avenue"
$newHotelSheet.Cells.Item(2, 8).Value = "1027"
$newHotelSheet.Cells.Item(2, 9).Value = "91"
$newHotelSheet.Cells.Item(2, 10).Value = "1047"
